$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.381.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6299'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07628'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07735'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.851.77'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.71%  '
$ws.Range('E13').Value = '  +12.33%  '
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6787'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.106.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.403.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.489'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.346'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.464'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.299'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05584'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.110'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.031'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.846'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.157'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7096'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.583'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.239.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01803'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.413'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.26%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.148'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4015'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.047'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.680'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1121'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.02%  '
